$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Redact / replace sensitive values in column B with sanitized placeholders ---
$ws.Range("B1").Value  = "xxxxxxxxxxxxxxxFe[4:@vil"
$ws.Range("B2").Value  = "https://xxxxxxx.crm8.dynamics.com or https://igdcicd2.api.crm.dynamics.com/XRMServices/2011/Organization.svc"
$ws.Range("B3").Value  = "xxxxxxxxxxxxxxx"
$ws.Range("B4").Value  = "xxxxxxx@IGDCRM.onmicrosoft.com"
$ws.Range("B6").Value  = "hxxxxxxxxxxxcf2hkkb4sl4tb2k37ciqdovsh7zq"
$ws.Range("B7").Value  = "xxxxxxx@IGDCRM.onmicrosoft.com"
$ws.Range("B8").Value  = "xxxxxxxxx-d509-4a1d-babf-xxxxxxxxx"
$ws.Range("B9").Value  = "xxxxxxxxx-a006-4359-966b-xxxxxxxxxxxx"
$ws.Range("B12").Value = "https://KeyVaultName.vault.azure.net/"

# --- Add new hyperlinks on GitUserName (B7) and BASESECRETURI (B12) ---
$ws.Hyperlinks.Add($ws.Range("B7"), "mailto:xxxxxxx@IGDCRM.onmicrosoft.com")
$ws.Hyperlinks.Add($ws.Range("B12"), "https://KeyVaultName.vault.azure.net/")

# Match the existing hyperlink look (same formatting as B4, the original hyperlink cell)
$ws.Range("B4").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)

# --- Update sheet view: scroll back to top-left, select B6 ---
$ws.Range("B6").Select()

Write-Host "done"
